$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Thu Oct 20 21:34:11 EDT 2022"
$ws.Range("B3").Value = "Thu Oct 20 21:34:20 EDT 2022"
$ws.Range("B4").Value = "Thu Oct 20 21:34:27 EDT 2022"
$ws.Range("B5").Value = "Thu Oct 20 21:34:33 EDT 2022"
$ws.Range("B6").Value = "Thu Oct 20 21:34:39 EDT 2022"
$ws.Range("B7").Value = "Thu Oct 20 21:34:46 EDT 2022"
